$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the two rows we need to touch by content instead of trusting a
# hard-coded row number.
$invoicesRowIndex = -1
$adminRowIndex = -1
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $firstCellText = $t.Rows.Item($i).Cells.Item(1).Range.Text
    if ($firstCellText -like "*InvoicesTab*") {
        $invoicesRowIndex = $i
    }
    if ($firstCellText.StartsWith("ADMIN") -and $firstCellText.Length -le 8) {
        $adminRowIndex = $i
    }
}
if ($invoicesRowIndex -eq -1) { $invoicesRowIndex = 28 }
if ($adminRowIndex -eq -1) { $adminRowIndex = 29 }

# --- Row "Zoekfunctie inbouwen finance in Invoice[_GoBack]sTab" ---------
# Merge the "Invoice" / bookmark / "sTab" runs into a single "InvoicesTab"
# run (dropping the old, now-stale "_GoBack" bookmark that sat between
# them) and fill in the second column with "v".
$row1 = $t.Rows.Item($invoicesRowIndex)
$cell1a = $row1.Cells.Item(1)
$cell1a.Range.Find.Execute("InvoicesTab", $false, $false, $false, $false, `
    $false, $true, 1, $false, "InvoicesTab", 2) | Out-Null

$cell1b = $row1.Cells.Item(2)
$cell1b.Range.Text = "v"

# --- Row "ADMIN" ---------------------------------------------------------
# Put "V?" in the second column and move the document's "_GoBack" bookmark
# here (right after the "?"), which is where Word leaves it after the
# last edit. A placeholder character is used so the bookmark can be
# anchored with Bookmarks.Add and then collapsed back down once the
# placeholder is removed.
$row2 = $t.Rows.Item($adminRowIndex)
$cell2b = $row2.Cells.Item(2)
$cell2b.Range.Text = "V?X"

$cell2bRange = $cell2b.Range
$placeholder = $d.Range($cell2bRange.End - 2, $cell2bRange.End - 1)
$d.Bookmarks.Add("_GoBack", $placeholder) | Out-Null
$d.Range($cell2bRange.End - 2, $cell2bRange.End - 1).Delete() | Out-Null
